# "Add files via upload" — add the two new "Lactose intolerant" remarks for
# the Lukas Larson row (row 12) in the roster table, and leave the
# selection where the author left it when they saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E12").Value = "I'm Lactose intolerant"
$ws.Range("F12").Value = "Still Lactose intolerant"

$ws.Range("E17").Select() | Out-Null
